$wb = $excel.ActiveWorkbook

# --- Update metadata on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "6.1.0"
$meta.Range("B8").Value = "2022-05-31T20:10:14+00:00"

# --- Remove the now-obsolete "Extension.extension.extension.*" rows (12-15) ---
# on the "Elements" sheet. This collapses an extra nesting level that was
# removed from the StructureDefinition, shifting old rows 16-21 up to 12-17.
$elements = $wb.Worksheets.Item("Elements")
$elements.Rows("12:15").Delete()
